$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has:
#   Col A: segment label (string), rows 2-20
#   Col B: PercActivationsOriginal
#   Col C: PercActivationsCorrect
#   Col D: PercActivationsIncorrect
#   Col E: PercActivationsCorrected
#   Col F: PercActivationsFixed
#
# Target layout adds a new "segments" index column at A (0-based row
# index) and moves the label text that used to live in column A into a
# new column B, shifting the old B:F data columns one slot to the right
# (C:G).

# Insert a new column before B - this shifts B:F -> C:G and also moves
# the header cells (and picks up left-column formatting on the new col,
# which we fix up below).
$ws.Columns("B:B").Insert()

# The inserted B column inherited column A's bold/border header style;
# the data rows (2-20) should have no special style (matching the old
# "value" columns), so strip that back off.
$ws.Range("B2:B20").ClearFormats()

# B1 needs the same header style as the other header cells - grab it
# from the neighboring header cell before writing the new title text.
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$ws.Range("B1").Value = "segments"

# New column A: 0-based segment index (numeric), keeps existing bold
# header style that previously belonged to the label column.
$indices = @(0,1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18)
for ($i = 0; $i -lt $indices.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $indices[$i]
}

# New column B: the segment label text that used to live in column A.
$labels = @(
    "background",
    "back_bumper",
    "back_glass",
    "back_left_door",
    "back_left_light",
    "back_right_door",
    "back_right_light",
    "front_bumper",
    "front_glass",
    "front_left_door",
    "front_left_light",
    "front_right_door",
    "front_right_light",
    "hood",
    "left_mirror",
    "right_mirror",
    "tailgate",
    "trunk",
    "wheel"
)
for ($i = 0; $i -lt $labels.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $labels[$i]
}
